$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wording: "implementacija sa User" -> "implementacija za User"
$ws.Range("A10").Value = "Kreiranje intefejsa aplikativnih servisnih metoda I implementacija za User"

# Fill in the actual time spent for the two completed User-related tasks
$ws.Range("C10").Value = "20min"
$ws.Range("C11").Value = "20min"

# Move selection to C10 (where the last edit was made)
$ws.Range("C10").Select()
